$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh rotated the content of rows 3, 4, 6, 7 and 8
# (row 2 and row 5 stay untouched). Apply the new values cell-by-cell so
# that each row ends up with the data described in the update.

# Row 3
$ws.Range("D3").Value = 44189
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16562
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("S3").Value = 920
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44181
$ws.Range("K4").Value = "Modesto"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1139

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("K6").Value = "Dina"
$ws.Range("M6").Value = 55
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15455
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("S6").Value = 1030
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44168
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 16500
$ws.Range("Q7").Value = '$/caja 16 kilos granel'
$ws.Range("S7").Value = 1031
$ws.Range("T7").Value = 16

# Row 8
$ws.Range("D8").Value = 44176
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17400
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 967
$ws.Range("T8").Value = 18

$wb.Save()
